# Apply ENTSOE Solar Production Historical data refresh
# - Shifts the two days of quarter-hourly data from 10-11 Apr 2025 to 28-29 Apr 2025 (+18 days)
# - Updates the "Lookup" column (E) text to match the new dates
# - Refreshes Notified/Actual Production (B/C) values with the newly retrieved readings

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2

# New Timestamp (col A) serial values, one per data row starting at row 2
$newTimestamps = @(
    "45775","45775.01041666666","45775.02083333334","45775.03125","45775.04166666666","45775.05208333334","45775.0625","45775.07291666666","45775.08333333334","45775.09375",
    "45775.10416666666","45775.11458333334","45775.125","45775.13541666666","45775.14583333334","45775.15625","45775.16666666666","45775.17708333334","45775.1875","45775.19791666666",
    "45775.20833333334","45775.21875","45775.22916666666","45775.23958333334","45775.25","45775.26041666666","45775.27083333334","45775.28125","45775.29166666666","45775.30208333334",
    "45775.3125","45775.32291666666","45775.33333333334","45775.34375","45775.35416666666","45775.36458333334","45775.375","45775.38541666666","45775.39583333334","45775.40625",
    "45775.41666666666","45775.42708333334","45775.4375","45775.44791666666","45775.45833333334","45775.46875","45775.47916666666","45775.48958333334","45775.5","45775.51041666666",
    "45775.52083333334","45775.53125","45775.54166666666","45775.55208333334","45775.5625","45775.57291666666","45775.58333333334","45775.59375","45775.60416666666","45775.61458333334",
    "45775.625","45775.63541666666","45775.64583333334","45775.65625","45775.66666666666","45775.67708333334","45775.6875","45775.69791666666","45775.70833333334","45775.71875",
    "45775.72916666666","45775.73958333334","45775.75","45775.76041666666","45775.77083333334","45775.78125","45775.79166666666","45775.80208333334","45775.8125","45775.82291666666",
    "45775.83333333334","45775.84375","45775.85416666666","45775.86458333334","45775.875","45775.88541666666","45775.89583333334","45775.90625","45775.91666666666","45775.92708333334",
    "45775.9375","45775.94791666666","45775.95833333334","45775.96875","45775.97916666666","45775.98958333334","45776","45776.01041666666","45776.02083333334","45776.03125",
    "45776.04166666666","45776.05208333334","45776.0625","45776.07291666666","45776.08333333334","45776.09375","45776.10416666666","45776.11458333334","45776.125","45776.13541666666",
    "45776.14583333334","45776.15625","45776.16666666666","45776.17708333334","45776.1875","45776.19791666666","45776.20833333334","45776.21875","45776.22916666666","45776.23958333334",
    "45776.25","45776.26041666666","45776.27083333334","45776.28125","45776.29166666666","45776.30208333334","45776.3125","45776.32291666666","45776.33333333334","45776.34375",
    "45776.35416666666","45776.36458333334","45776.375","45776.38541666666","45776.39583333334","45776.40625","45776.41666666666","45776.42708333334","45776.4375","45776.44791666666",
    "45776.45833333334","45776.46875","45776.47916666666","45776.48958333334","45776.5","45776.51041666666","45776.52083333334","45776.53125","45776.54166666666","45776.55208333334",
    "45776.5625","45776.57291666666","45776.58333333334","45776.59375","45776.60416666666","45776.61458333334","45776.625","45776.63541666666","45776.64583333334","45776.65625",
    "45776.66666666666","45776.67708333334","45776.6875","45776.69791666666","45776.70833333334","45776.71875","45776.72916666666","45776.73958333334","45776.75","45776.76041666666",
    "45776.77083333334","45776.78125","45776.79166666666","45776.80208333334","45776.8125","45776.82291666666","45776.83333333334","45776.84375","45776.85416666666","45776.86458333334",
    "45776.875","45776.88541666666","45776.89583333334","45776.90625","45776.91666666666","45776.92708333334","45776.9375","45776.94791666666","45776.95833333334","45776.96875",
    "45776.97916666666","45776.98958333334"
)

# New Lookup (col E) text values, one per data row starting at row 2
$newLookups = @(
    "28.04.20251","28.04.20252","28.04.20253","28.04.20254","28.04.20255","28.04.20256","28.04.20257","28.04.20258","28.04.20259","28.04.202510",
    "28.04.202511","28.04.202512","28.04.202513","28.04.202514","28.04.202515","28.04.202516","28.04.202517","28.04.202518","28.04.202519","28.04.202520",
    "28.04.202521","28.04.202522","28.04.202523","28.04.202524","28.04.202525","28.04.202526","28.04.202527","28.04.202528","28.04.202529","28.04.202530",
    "28.04.202531","28.04.202532","28.04.202533","28.04.202534","28.04.202535","28.04.202536","28.04.202537","28.04.202538","28.04.202539","28.04.202540",
    "28.04.202541","28.04.202542","28.04.202543","28.04.202544","28.04.202545","28.04.202546","28.04.202547","28.04.202548","28.04.202549","28.04.202550",
    "28.04.202551","28.04.202552","28.04.202553","28.04.202554","28.04.202555","28.04.202556","28.04.202557","28.04.202558","28.04.202559","28.04.202560",
    "28.04.202561","28.04.202562","28.04.202563","28.04.202564","28.04.202565","28.04.202566","28.04.202567","28.04.202568","28.04.202569","28.04.202570",
    "28.04.202571","28.04.202572","28.04.202573","28.04.202574","28.04.202575","28.04.202576","28.04.202577","28.04.202578","28.04.202579","28.04.202580",
    "28.04.202581","28.04.202582","28.04.202583","28.04.202584","28.04.202585","28.04.202586","28.04.202587","28.04.202588","28.04.202589","28.04.202590",
    "28.04.202591","28.04.202592","28.04.202593","28.04.202594","28.04.202595","28.04.202596","29.04.20251","29.04.20252","29.04.20253","29.04.20254",
    "29.04.20255","29.04.20256","29.04.20257","29.04.20258","29.04.20259","29.04.202510","29.04.202511","29.04.202512","29.04.202513","29.04.202514",
    "29.04.202515","29.04.202516","29.04.202517","29.04.202518","29.04.202519","29.04.202520","29.04.202521","29.04.202522","29.04.202523","29.04.202524",
    "29.04.202525","29.04.202526","29.04.202527","29.04.202528","29.04.202529","29.04.202530","29.04.202531","29.04.202532","29.04.202533","29.04.202534",
    "29.04.202535","29.04.202536","29.04.202537","29.04.202538","29.04.202539","29.04.202540","29.04.202541","29.04.202542","29.04.202543","29.04.202544",
    "29.04.202545","29.04.202546","29.04.202547","29.04.202548","29.04.202549","29.04.202550","29.04.202551","29.04.202552","29.04.202553","29.04.202554",
    "29.04.202555","29.04.202556","29.04.202557","29.04.202558","29.04.202559","29.04.202560","29.04.202561","29.04.202562","29.04.202563","29.04.202564",
    "29.04.202565","29.04.202566","29.04.202567","29.04.202568","29.04.202569","29.04.202570","29.04.202571","29.04.202572","29.04.202573","29.04.202574",
    "29.04.202575","29.04.202576","29.04.202577","29.04.202578","29.04.202579","29.04.202580","29.04.202581","29.04.202582","29.04.202583","29.04.202584",
    "29.04.202585","29.04.202586","29.04.202587","29.04.202588","29.04.202589","29.04.202590","29.04.202591","29.04.202592","29.04.202593","29.04.202594",
    "29.04.202595","29.04.202596"
)

for ($i = 0; $i -lt $newTimestamps.Length; $i++) {
    $row = $firstDataRow + $i
    $ws.Cells.Item($row, 1).Value2 = [double]$newTimestamps[$i]
    $ws.Cells.Item($row, 5).Value = $newLookups[$i]
}

# Refreshed Notified Production (col B) / Actual Production (col C) readings.
# Each entry is (row, newNotifiedValue, newActualValue); a value of $null means "unchanged".
$bcChanges = @(
    @(22,20,$null),
    @(23,21,2),
    @(24,26,15),
    @(25,38,37),
    @(26,254,94),
    @(27,281,150),
    @(28,311,230),
    @(29,345,308),
    @(30,744,452),
    @(31,803,550),
    @(32,853,618),
    @(33,914,691),
    @(34,1393,880),
    @(35,1430,997),
    @(36,1472,1068),
    @(37,1515,1095),
    @(38,1830,1274),
    @(39,1865,1305),
    @(40,1894,1325),
    @(41,1925,1367),
    @(42,2129,1526),
    @(43,2150,1595),
    @(44,2169,1587),
    @(45,2185,1623),
    @(46,1973,1381),
    @(47,1981,1310),
    @(48,1987,1331),
    @(49,1990,1340),
    @(50,1973,1334),
    @(51,1973,1361),
    @(52,1971,1341),
    @(53,1967,1330),
    @(54,1919,1293),
    @(55,1911,1292),
    @(56,1902,1232),
    @(57,1890,1142),
    @(58,1781,1173),
    @(59,1763,1190),
    @(60,1744,1215),
    @(61,1719,1160),
    @(62,1531,1050),
    @(63,1498,1029),
    @(64,1467,999),
    @(65,1434,937),
    @(66,1305,1071),
    @(67,1265,1051),
    @(68,1231,950),
    @(69,1194,854),
    @(70,698,690),
    @(71,648,568),
    @(72,596,465),
    @(73,550,364),
    @(74,182,224),
    @(75,149,140),
    @(76,123,82),
    @(77,104,41),
    @(78,26,16),
    @(79,21,3),
    @(80,20,13),
    @(81,20,7),
    @(82,17,3),
    @(83,17,$null),
    @(84,17,$null),
    @(85,17,$null),
    @(86,17,$null),
    @(87,16,$null),
    @(88,16,$null),
    @(89,16,$null),
    @(118,18,$null),
    @(119,19,4),
    @(120,25,23),
    @(121,37,49),
    @(122,287,107),
    @(123,316,172),
    @(124,351,249),
    @(125,392,348),
    @(126,967,513),
    @(127,1016,619),
    @(128,1069,717),
    @(129,1126,814),
    @(130,1583,984),
    @(131,1624,1095),
    @(132,1673,1157),
    @(133,1714,1169),
    @(134,1995,1294),
    @(135,2026,$null),
    @(136,2056,$null),
    @(137,2082,$null),
    @(138,2257,$null),
    @(139,2273,$null),
    @(140,2289,$null),
    @(141,2302,$null),
    @(142,2077,$null),
    @(143,2083,$null),
    @(144,2086,$null),
    @(145,2087,$null),
    @(146,2087,$null),
    @(147,2086,$null),
    @(148,2081,$null),
    @(149,2077,$null),
    @(150,1998,$null),
    @(151,1988,$null),
    @(152,1977,$null),
    @(153,1963,$null),
    @(154,1862,$null),
    @(155,1842,$null),
    @(156,1817,$null),
    @(157,1796,$null),
    @(158,1687,$null),
    @(159,1651,$null),
    @(160,1619,$null),
    @(161,1582,$null),
    @(162,1411,$null),
    @(163,1369,$null),
    @(164,1334,$null),
    @(165,1295,$null),
    @(166,810,$null),
    @(167,756,$null),
    @(168,700,$null),
    @(169,660,$null),
    @(170,255,$null),
    @(171,224,$null),
    @(172,202,$null),
    @(173,187,$null),
    @(174,27,$null),
    @(175,21,$null),
    @(176,21,$null),
    @(177,20,$null),
    @(178,17,$null),
    @(179,17,$null),
    @(180,17,$null),
    @(181,17,$null),
    @(182,17,$null),
    @(183,16,$null),
    @(184,16,$null),
    @(185,16,$null)
)

foreach ($entry in $bcChanges) {
    $row = $entry[0]
    $newB = $entry[1]
    $newC = $entry[2]
    if ($null -ne $newB) { $ws.Cells.Item($row, 2).Value = $newB }
    if ($null -ne $newC) { $ws.Cells.Item($row, 3).Value = $newC }
}

Write-Host "Applied Solar_Production_Historical refresh: shifted dates +18 days and updated production readings."